$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 46, shifting existing rows 46:110 down to 47:111
$ws.Rows("46:46").Insert()

# Populate the new row 46 with the new data record (mirrors the other Cilantro rows' static columns)
$ws.Range("A46").Value = 8
$ws.Range("B46").Value = "Terminal La Palmera de La Serena"
$ws.Range("C46").Value = "Coquimbo"
$ws.Range("D46").Value = 44579
$ws.Range("E46").Value = 4
$ws.Range("F46").Value = 100112040
$ws.Range("G46").Value = "Cilantro"
$ws.Range("H46").Value = "Sin especificar"
$ws.Range("I46").Value = "Primera"
$ws.Range("J46").Value = 3360
$ws.Range("K46").Value = 2500
$ws.Range("L46").Value = 3000
$ws.Range("M46").Value = 2750
$ws.Range("N46").Value = "`$/atado 1 a 1,5 kilos"
$ws.Range("O46").Value = "Provincia del Elquí"
$ws.Range("P46").Value = 1833
$ws.Range("Q46").Value = 1.5
$ws.Range("R46").Value = "Hortaliza"
